$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new account row ---
# New row goes immediately above the existing "004222784" (RAFAEL) row,
# i.e. right after "004580355" (LARISSA) in the Conta-sorted export.
$anchor = $ws.Columns.Item(1).Find("004222784")
$newRow = $anchor.Row

$ws.Rows.Item($newRow).Insert()

# Format the account-number cell as text first so the leading zeros in
# "005624730" are preserved instead of the value being parsed as a number.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "005624730"
$ws.Cells.Item($newRow, 2).Value = "ISABEL"
$ws.Cells.Item($newRow, 3).Value = 20000

# --- Remove the account row that was dropped from the export ---
# Account 004237325 / RICARDO / -11636.77
$toRemove = $ws.Columns.Item(1).Find("004237325")
$ws.Rows.Item($toRemove.Row).Delete()
